$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数 / interested-count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1095
$ws1.Range("F3").Value = 4117
$ws1.Range("F7").Value = 14
$ws1.Range("F8").Value = 32
$ws1.Range("F10").Value = 122
$ws1.Range("F11").Value = 299
$ws1.Range("F12").Value = 228
$ws1.Range("F13").Value = 2883
$ws1.Range("F14").Value = 128
$ws1.Range("F15").Value = 1365
$ws1.Range("F16").Value = 9

# Sheet "全部类型" (All Types) - update column F (想去人数 / interested-count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1095
$ws4.Range("F3").Value = 4117
$ws4.Range("F8").Value = 14
$ws4.Range("F9").Value = 32
$ws4.Range("F11").Value = 122
$ws4.Range("F12").Value = 299
$ws4.Range("F13").Value = 228
$ws4.Range("F14").Value = 2883
$ws4.Range("F15").Value = 128
$ws4.Range("F16").Value = 1365
$ws4.Range("F17").Value = 9
